$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E14").Value = $null

$ws.Range("B14").Value = 234234
$ws.Range("D14").Value = 234234

$ws.Range("C16").Value = "qweq3"

$ws.Range("B20").Value = 234234
$ws.Range("E20").Value = "433t34t34t"

$ws.Range("E20").Select()
